$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(116).Insert()

$ws.Range("A116").Value = 4
$ws.Range("B116").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C116").Value = "Los Lagos"
$ws.Range("D116").Value = 44582
$ws.Range("E116").Value = 10
$ws.Range("F116").Value = "Fruta"
$ws.Range("G116").Value = 100108
$ws.Range("H116").Value = "Tropicales y subtropicales"
$ws.Range("I116").Value = 100108005
$ws.Range("J116").Value = "Piña"
$ws.Range("K116").Value = "Caramelo"
$ws.Range("L116").Value = "Tercera"
$ws.Range("M116").Value = 200
$ws.Range("N116").Value = 18000
$ws.Range("O116").Value = 19000
$ws.Range("P116").Value = 18500
$ws.Range("Q116").Value = "$/caja 16 unidades"
$ws.Range("R116").Value = "Ecuador"
$ws.Range("S116").Value = 1156
$ws.Range("T116").Value = 16
